# Apply the "worknotes" update:
#  - TODO sheet: add a new follow-up note row about whether store user info
#    survives a page reload after login (tracks removal of org_dept_level
#    from the user db).
#  - 机组启停记录 (unit start/stop log) sheet: rename the "停止标志" status
#    column to "运行标志"/isRunning, drop the obsolete "2-已停止" legend row,
#    and append a new flow description (rows 25-33) for the
#    getGenStartLog API used to look up a unit's running state.
#  - Refresh selections / active sheet to match the saved workbook state.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "TODO"
# ---------------------------------------------------------------------
$todo = $wb.Worksheets.Item("TODO")

$todo.Range("A4").Value = 3
$todo.Range("B4").Value = 43578
$todo.Range("B2").Copy()
$todo.Range("B4").PasteSpecial(-4122)
$todo.Range("C4").Value = "用户登录成功，未点击退出关闭浏览器或窗口，Token未过期，再打开网页，store的用户信息是否丢失"

# ---------------------------------------------------------------------
# Sheet "机组启停记录" (Unit start/stop log)
# ---------------------------------------------------------------------
$log = $wb.Worksheets.Item("机组启停记录")

# Relabel the "stop flag" legend as a "running flag" legend and drop the
# now-unused "2-已停止" row.
$log.Range("G5").Value = "运行标志"
$log.Range("G6").Value = "isRunning"
$log.Range("G7").ClearContents()

# New flow block describing the getGenStartLog lookup.
$log.Range("C25").Value = "流程"

$log.Range("A26").Value = 43578
$log.Range("A1").Copy()
$log.Range("A26").PasteSpecial(-4122)
$log.Range("C26").Value = "page-1"

$log.Range("D27").Value = "页面加载，识别用户，获取stationIdx"
$log.Range("H27").Value = "controller - Devices"

$log.Range("I28").Value = "method -get_gen_start_log"

$log.Range("D29").Value = "选择机组编号，向server查找记录"
$log.Range("I29").Value = "查找最后一条记录"
$log.Range("I29").Characters(3, 4).Font.Color = 255

$log.Range("E30").Value = "api-getGenStartLog"
$log.Range("J30").Value = "isRunning = TRUE，机组运行，回应启动时间（不可编辑），要求填写停止时间"

$log.Range("F31").Value = "para"
$log.Range("J31").Value = "isRunning = FALSE，机组停止，填写启动时间"

$log.Range("F32").Value = "stationIdx"

$log.Range("F33").Value = "genIdx"

# ---------------------------------------------------------------------
# Restore view state: selections per sheet, and the active sheet/tab.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("注册").Select()
$wb.Worksheets.Item("注册").Range("G15").Select()

$wb.Worksheets.Item("登录").Select()
$wb.Worksheets.Item("登录").Range("D13").Select()

$wb.Worksheets.Item("用户设置").Select()
$wb.Worksheets.Item("用户设置").Range("G22").Select()

$log.Select()
$log.Range("D27").Select()

$wb.Worksheets.Item("全局").Select()
$wb.Worksheets.Item("全局").Range("H23").Select()

$todo.Select()
$todo.Range("F10").Select()

$store = $wb.Worksheets.Item("Store")
$store.Select()
$store.Range("E12").Select()
